# Implementation for relationship classes
# Adds a new worksheet "DatapointObservesSpatialElement" (the
# Datapoint-observes-SpatialElement relationship class table) after the
# existing PressureDatapoint sheet, fills it with its header + two rows,
# and nudges the selections on the affected sheets to match.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "DatapointObservesSpatialElement"

# Header row.
$newSheet.Range("A1").Value = "name"
$newSheet.Range("B1").Value = "deviceid"

# Data rows - relates the two TemperatureDatapoint rows (T-1, T-2) to
# Device id 4.
$newSheet.Range("A2").Value = "T-1"
$newSheet.Range("B2").Value = 4
$newSheet.Range("A3").Value = "T-2"
$newSheet.Range("B3").Value = 4

$newSheet.PageSetup.Orientation = 1

# TemperatureDatapoint's selection moves from B2 to A2.
$tempSheet = $wb.Worksheets.Item("TemperatureDatapoint")
$tempSheet.Activate() | Out-Null
$tempSheet.Range("A2").Select() | Out-Null

# Leave the new sheet active/selected (tab 4, cell A3) - matches the
# workbook's saved view state.
$newSheet.Activate() | Out-Null
$newSheet.Range("A3").Select() | Out-Null
